$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update "想去人数" (F column) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 13395
$ws1.Range("F11").Value = 70
$ws1.Range("F14").Value = 13371
$ws1.Range("F16").Value = 588
$ws1.Range("F17").Value = 8913
$ws1.Range("F19").Value = 7983

# Sheet "全部类型" (sheet4): same underlying rows, offset by one row
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 13395
$ws4.Range("F12").Value = 70
$ws4.Range("F15").Value = 13371
$ws4.Range("F17").Value = 588
$ws4.Range("F18").Value = 8913
$ws4.Range("F20").Value = 7983
